# Update "想去人数" (number of people wishing to attend) figures across the
# "展览" and "全部类型" worksheets to reflect freshly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 263
$ws1.Range("F5").Value = 6634
$ws1.Range("F6").Value = 5416
$ws1.Range("F8").Value = 67
$ws1.Range("F11").Value = 234
$ws1.Range("F12").Value = 48

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 263
$ws4.Range("F5").Value = 6634
$ws4.Range("F6").Value = 5416
$ws4.Range("F8").Value = 67
$ws4.Range("F11").Value = 234
$ws4.Range("F14").Value = 48
